$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove the "H 72" record (row 2); all subsequent rows shift up by one.
$ws.Rows.Item(2).Delete()

# Step 2: apply the updated/re-sampled missing-data pattern for the remaining rows.
$ws.Range("E2").Value = -6.9
$ws.Range("D4").Value = -13.5
$ws.Range("E5").Value = ""
$ws.Range("E7").Value = -5
$ws.Range("D8").Value = ""
$ws.Range("D9").Value = -15.4
$ws.Range("B10").Value = -19.5
$ws.Range("D10").Value = -13.8
$ws.Range("E10").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = -6.6
$ws.Range("B12").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E13").Value = -6.1
$ws.Range("B14").Value = -20.8
$ws.Range("E14").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("E15").Value = -12
$ws.Range("E18").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("D20").Value = -15.3
$ws.Range("E21").Value = -7.3
$ws.Range("D22").Value = ""
$ws.Range("D23").Value = -15.5
$ws.Range("E24").Value = ""
$ws.Range("D25").Value = ""
$ws.Range("B26").Value = -19.5
$ws.Range("B28").Value = ""
$ws.Range("D28").Value = -13.9
$ws.Range("E28").Value = -7
$ws.Range("B30").Value = -19.5
$ws.Range("E31").Value = ""
$ws.Range("B32").Value = ""
$ws.Range("D32").Value = ""
$ws.Range("E33").Value = -9.5
$ws.Range("B35").Value = -19.2
$ws.Range("B36").Value = ""
$ws.Range("D36").Value = -14.6
$ws.Range("E36").Value = ""
$ws.Range("B37").Value = -19.8
$ws.Range("B38").Value = ""
$ws.Range("D38").Value = ""
$ws.Range("D41").Value = -15.9
$ws.Range("D42").Value = -14.6
$ws.Range("D43").Value = ""
$ws.Range("D44").Value = ""
$ws.Range("B45").Value = -19.7
$ws.Range("B46").Value = ""
$ws.Range("D52").Value = -13.8
$ws.Range("B53").Value = -20.3
$ws.Range("D54").Value = ""
$ws.Range("B56").Value = ""
$ws.Range("E57").Value = -5.9
$ws.Range("E59").Value = -5.7
$ws.Range("E60").Value = ""
$ws.Range("E62").Value = ""
